$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'290.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.20%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'29.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'4.12%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.130"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'4.33%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06685"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.04%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.360"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.53%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.353"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-1.70%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9153"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.11%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1590"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.09%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.06671"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.17%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07674"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.11%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.02942"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.28%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.08997"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.21%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001590"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.04%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04510"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.02%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").Value = "'0.0006461"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.17%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006257"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'3.79%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.450"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.61%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").Value = "'3.406"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.33%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-0.91%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D23").Value = "'4.067"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'2.43%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D25").Value = "'0.001193"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.52%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004133"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-4.74%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'5.89%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'-1.10%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.04231"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.69%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006719"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.43%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'-12.09%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-3.91%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01161"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-6.51%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005627"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.60%"
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "BOLO"
$ws.Range("C46").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D46").Value = "'1.974"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'26.36%"
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "'0.01306"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-29.40%"
$ws.Range("E47").Style = "Normal"
